# dataText.xlsx fixes:
#  - tidy up the metric-name labels in column A (remove stray leading
#    double-space / trailing-space text, use camel-case tokens instead)
#  - scroll the sheet down a bit and leave the cursor on E6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "fix white space in buckinghamshire" / "fix data text metric names tidycode":
# replace the three untidy, whitespace-padded labels with clean identifiers.
$ws.Range("A7").Value = "SelfEmployed"
$ws.Range("A8").Value = "Unemployed"
$ws.Range("A9").Value = "Inactive"

# Scroll the view down so row 5 is at the top and select E6, matching the
# author's saved window position.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("E6").Select()
